$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string used by B1 header ("motivation" -> "use")
$ws.Range("B1").Value = "use"

# Update the data values (randomization was added to study)
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 3

$ws.Range("A3").Value = 7
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 4

$ws.Range("A4").Value = 8
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 3

$ws.Range("A5").Value = 9
$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 6
$ws.Range("C6").Value = 4

# Update the selected cell in the sheet view
$ws.Range("B2").Select()
